# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# - Updates "VALOR MORA" (E11) and "Cant. Periodos" (F13).
# - Inserts a new worker-period row (period 2508) below the existing
#   2507/2506/2505/2504 rows, re-sorting the period rows into ascending
#   order (2504..2508) and updating the "Salario Basico" (column G)
#   amount for every period row.
# - The signature block (old rows 24-25) naturally shifts down to rows
#   25-26 because of the row insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header figures -------------------------------------------------
$ws.Range("E11").Value = 284700
$ws.Range("F13").Value = 5

# --- Make room for the new period row (2508) -------------------------
# Inserting at row 20 pushes the signature block (rows 24-25) down to
# rows 25-26 and leaves a blank row 20 ready for the new data.
$ws.Rows("20").Insert()

# Row 20 should carry the "last row" border styling that row 19
# currently has (thicker outer border), so copy its formatting down
# before row 19 gets reformatted as a normal middle row.
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)

# Row 19 becomes a normal middle row, matching rows 16-18.
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)

# --- Re-sort the period rows into ascending order and bump salary ---
$ws.Range("E16").Value = "2504"
$ws.Range("E17").Value = "2505"
$ws.Range("E18").Value = "2506"
$ws.Range("E19").Value = "2507"

$ws.Range("G16:G19").Value = 1423500

# --- New row 20: period 2508 -----------------------------------------
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1050966338"
$ws.Range("D20").Value = "MARCIO ENRIQUE HUETO ZAMBRANO"
$ws.Range("E20").Value = "2508"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500
